$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.141.92'
$ws.Range("E2").Value = '  -5.05%  '
$ws.Range("D3").Value = '3.702.54'
$ws.Range("E3").Value = '  -4.70%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '583.66'
$ws.Range("E5").Value = '  -2.36%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '180.18'
$ws.Range("E6").Value = '  +7.58%  '
$ws.Range("D7").Value = '3.692.08'
$ws.Range("E7").Value = '  -4.87%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.628'
$ws.Range("E8").Value = '  -6.08%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.00'
$ws.Range("E9").Value = '  +0.30%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.713'
$ws.Range("E10").Value = '  -6.23%  '
$ws.Range("E11").Value = '  -8.60%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '53.87'
$ws.Range("E12").Value = '  -0.23%  '
$ws.Range("E13").Value = '  -9.75%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '10.43'
$ws.Range("E14").Value = '  -8.04%  '
$ws.Range("D15").Value = '4.355.64'
$ws.Range("E15").Value = '  -3.42%  '
$ws.Range("D16").Value = '3.734.00'
$ws.Range("E16").Value = '  -3.86%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '19.50'
$ws.Range("E17").Value = '  -7.12%  '
$ws.Range("E18").Value = '  -2.67%  '
$ws.Range("E19").Value = '  -7.69%  '
$ws.Range("E20").Value = '  -7.58%  '
$ws.Range("D21").Value = '67.876.37'
$ws.Range("E21").Value = '  -5.27%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '408.72'
$ws.Range("E22").Value = '  -6.10%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.50'
$ws.Range("E23").Value = '  -4.91%  '
$ws.Range("E24").Value = '  -6.07%  '
$ws.Range("E25").Value = '  -8.29%  '
$ws.Range("E26").Value = '  -7.78%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.00'
$ws.Range("E27").Value = '  +0.26%  '
$ws.Range("E28").Value = '  -7.67%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.07'
$ws.Range("E29").Value = '  +2.45%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.50'
$ws.Range("E30").Value = '  -6.79%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '32.54'
$ws.Range("E31").Value = '  -7.27%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.45'
$ws.Range("E32").Value = '  -7.39%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '12.51'
$ws.Range("E33").Value = '  -8.17%  '
$ws.Range("E34").Value = '  -7.12%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '65.25'
$ws.Range("E35").Value = '  -4.69%  '
$ws.Range("B36").Value = 'Bittensor'
$ws.Range("C36").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '600.95'
$ws.Range("E36").Value = '  -3.51%  '
$ws.Range("B37").Value = 'InjectiveProtocol'
$ws.Range("C37").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '42.91'
$ws.Range("E37").Value = '  -17.68%  '
$ws.Range("D38").Value = '0.0₃0895'
$ws.Range("E38").Value = '  -8.91%  '
$ws.Range("E39").Value = '  +0.17%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.397'
$ws.Range("E40").Value = '  -5.31%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.00'
$ws.Range("E41").Value = '  +0.05%  '
$ws.Range("E42").Value = '  -4.32%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.76'
$ws.Range("E43").Value = '  +5.07%  '
$ws.Range("E44").Value = '  -9.61%  '
$ws.Range("E45").Value = '  -8.65%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0435'
$ws.Range("E46").Value = '  -7.42%  '
$ws.Range("E47").Value = '  -9.77%  '
$ws.Range("D48").Value = '2.799.85'
$ws.Range("E48").Value = '  -2.41%  '
$ws.Range("E49").Value = '  -7.33%  '
$ws.Range("E50").Value = '  -4.76%  '
$ws.Range("E51").Value = '  -6.94%  '
